# Banking Project Brief.docx - apply the commit's edit.
#
# The diff centers the paragraph that holds the inline "Picture 3" image
# (the last paragraph in the document body): it adds <w:jc w:val="center"/>
# to that paragraph's <w:pPr>, right after the existing <w:ind w:hanging="11"/>.
#
# (The same diff also shows customXml/item2.xml <-> customXml/item4.xml and
# itemProps2.xml <-> itemProps4.xml swapping places - those are raw OPC
# package parts [a SharePoint "FormTemplates" part and a bibliography
# "b:Sources" part] that are not exposed anywhere on the Word object model
# -- Document.CustomXMLParts.Count is always 0 in this host -- so there is
# no COM-interop surface capable of touching them; only the paragraph
# formatting change below is reachable through Word automation.)

$d = $word.ActiveDocument

# Locate the paragraph that contains the document's inline picture. This is
# the final paragraph in the body (w14:paraId="0FAF56BB"), which currently
# has <w:pPr><w:ind w:hanging="11"/><w:rPr>...</w:rPr></w:pPr>.
if ($d.InlineShapes.Count -ge 1) {
    $targetParagraph = $d.InlineShapes.Item(1).Range.Paragraphs.First
} else {
    $targetParagraph = $d.Paragraphs.Last
}

# wdAlignParagraphCenter = 1
$targetParagraph.Alignment = 1

Write-Output ("Centered paragraph; Alignment=" + $targetParagraph.Alignment)
